$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "04. " prefix from the glossary title, leaving "Glossário".
# ---------------------------------------------------------------------------
$titleRange = $d.Range(0, 0)
$titleRange.Find.Execute("04. ") | Out-Null
$titleRange.Delete()

# ---------------------------------------------------------------------------
# 2. Remove the "Layout" and "Sketch" glossary rows entirely.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)

function Get-CellPlainText($cell) {
    $rawCellText = $cell.Range.Text
    return $rawCellText.TrimEnd([char]13, [char]7)
}

$rowsToRemove = @("Layout", "Sketch")
foreach ($termToRemove in $rowsToRemove) {
    for ($rowIdx = 1; $rowIdx -le $tbl.Rows.Count; $rowIdx++) {
        $candidateRow = $tbl.Rows.Item($rowIdx)
        $candidateText = Get-CellPlainText $candidateRow.Cells.Item(1)
        if ($candidateText -eq $termToRemove) {
            $candidateRow.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Update the "Paleta de Cores" definition text.
# ---------------------------------------------------------------------------
for ($rowIdx = 1; $rowIdx -le $tbl.Rows.Count; $rowIdx++) {
    $candidateRow = $tbl.Rows.Item($rowIdx)
    $candidateText = Get-CellPlainText $candidateRow.Cells.Item(1)
    if ($candidateText -eq "Paleta de Cores") {
        $defCell = $candidateRow.Cells.Item(2)
        $defCell.Range.Find.Execute("Em design gráfico, a paleta de cores", $true, $false, $false, $false, $false,
                                     $true, 0, $false, "A paleta de cores", 1) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 4. Update the "Suporte" definition text.
# ---------------------------------------------------------------------------
for ($rowIdx = 1; $rowIdx -le $tbl.Rows.Count; $rowIdx++) {
    $candidateRow = $tbl.Rows.Item($rowIdx)
    $candidateText = Get-CellPlainText $candidateRow.Cells.Item(1)
    if ($candidateText -eq "Suporte") {
        $defCell = $candidateRow.Cells.Item(2)
        $defCell.Range.Find.Execute("papel, tecido, etc.", $true, $false, $false, $false, $false,
                                     $true, 0, $false, "papel e tecido.", 1) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 5. Italicize "Mockup" (term cell) and "mockup" (inside the definition).
# ---------------------------------------------------------------------------
for ($rowIdx = 1; $rowIdx -le $tbl.Rows.Count; $rowIdx++) {
    $candidateRow = $tbl.Rows.Item($rowIdx)
    $candidateText = Get-CellPlainText $candidateRow.Cells.Item(1)
    if ($candidateText -eq "Mockup") {
        $termCell = $candidateRow.Cells.Item(1)
        $defCell = $candidateRow.Cells.Item(2)

        $termFind = $termCell.Range.Find
        $termFind.ClearFormatting()
        $termFind.Replacement.ClearFormatting()
        $termFind.Replacement.Font.Italic = $true
        $termFind.Execute("Mockup", $true, $false, $false, $false, $false,
                           $true, 0, $false, "Mockup", 1) | Out-Null

        $defFind = $defCell.Range.Find
        $defFind.ClearFormatting()
        $defFind.Replacement.ClearFormatting()
        $defFind.Replacement.Font.Italic = $true
        $defFind.Execute("mockup", $true, $false, $false, $false, $false,
                          $true, 0, $false, "mockup", 1) | Out-Null
        break
    }
}
